$wb = $excel.ActiveWorkbook

# Sheet handles (1-based index, matches workbook.xml <sheets> order)
$wsTransfer = $wb.Worksheets.Item(12)   # Multiple_FD_Transfer
$wsCash     = $wb.Worksheets.Item(13)   # Multiple_FD_Cash
$wsClosure  = $wb.Worksheets.Item(14)   # Closure_and_Renewals_Cash

# ---------------------------------------------------------------------------
# Multiple_FD_Cash (sheet13): drop the "LoginNewUserName"/"LoginNewPassword"
# columns (R:S) including the mailto hyperlink that lived in S2.
# ---------------------------------------------------------------------------
$wsCash.Hyperlinks.Delete()
$wsCash.Range("R1:S2").Clear()

# ---------------------------------------------------------------------------
# Closure_and_Renewals_Cash (sheet14): repurpose the sheet for the new
# "Jewel Gold Loan Renewal" data set - new headers, new sample row, drop the
# old hyperlink/extra columns, and widen the columns that now hold data.
# ---------------------------------------------------------------------------
$wsClosure.Hyperlinks.Delete()

$wsClosure.Range("E1").Value = "Item_Description"
$wsClosure.Range("F1").Value = "Item_Quantity"
$wsClosure.Range("G1").Value = "Stone_Weight"
$wsClosure.Range("H1").Value = "Item_Weight"
$wsClosure.Range("I1").Value = "Dirt_Weight"
$wsClosure.Range("J1").Value = "Enter_Purity"
$wsClosure.Range("K1:S1").ClearContents()

$wsClosure.Range("A2").Value = "Loan_Closure_Cash"
$wsClosure.Range("E2").Value = "necklace"
$wsClosure.Range("F2").Value = 25
$wsClosure.Range("G2").Value = 0
$wsClosure.Range("H2").Value = 24
$wsClosure.Range("I2").Value = 1
$wsClosure.Range("J2").Value = 22
$wsClosure.Range("K2:S2").ClearContents()

$wsClosure.Columns.Item(3).ColumnWidth = 22.85546875
$wsClosure.Columns.Item(4).ColumnWidth = 25.5703125
$wsClosure.Columns.Item(5).ColumnWidth = 29.140625
$wsClosure.Columns.Item(6).ColumnWidth = 20.7109375
$wsClosure.Columns.Item(7).ColumnWidth = 17.42578125
$wsClosure.Columns.Item(8).ColumnWidth = 17
$wsClosure.Columns.Item(9).ColumnWidth = 18.7109375
$wsClosure.Columns.Item(10).ColumnWidth = 13.42578125

# ---------------------------------------------------------------------------
# Selections / active sheet. The workbook now opens on
# Closure_and_Renewals_Cash (tab index 13, 0-based) instead of
# Multiple_FD_Transfer; Multiple_FD_Cash and Multiple_FD_Transfer keep their
# own remembered selections but are no longer the active tab.
# ---------------------------------------------------------------------------
$wsTransfer.Range("A10").Select()
$wsCash.Range("Q13").Select()
$wsClosure.Range("J8").Select()
